$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44489
$ws.Range("J2").Value = 600
$ws.Range("M2").Value = 1450
$ws.Range("P2").Value = 1450

# Row 4
$ws.Range("D4").Value = 44545
$ws.Range("K4").Value = 1700
$ws.Range("L4").Value = 1800
$ws.Range("M4").Value = 1755
$ws.Range("N4").Value = "$/kilo"
$ws.Range("P4").Value = 1755

# Row 5
$ws.Range("D5").Value = 44511
$ws.Range("J5").Value = 600
$ws.Range("K5").Value = 1300
$ws.Range("L5").Value = 1400
$ws.Range("M5").Value = 1350
$ws.Range("P5").Value = 1350

# Row 6
$ws.Range("D6").Value = 44510
$ws.Range("J6").Value = 600
$ws.Range("K6").Value = 1300
$ws.Range("L6").Value = 1400
$ws.Range("M6").Value = 1350
$ws.Range("P6").Value = 1350

# Row 7
$ws.Range("D7").Value = 44524
$ws.Range("J7").Value = 200
$ws.Range("O7").Value = "Provincia de Talca"

# Row 8
$ws.Range("D8").Value = 44519
$ws.Range("J8").Value = 250
$ws.Range("K8").Value = 1200
$ws.Range("L8").Value = 1300
$ws.Range("M8").Value = 1240
$ws.Range("P8").Value = 1240

# Row 9
$ws.Range("D9").Value = 44477
$ws.Range("J9").Value = 500
$ws.Range("K9").Value = 1400
$ws.Range("L9").Value = 1500
$ws.Range("M9").Value = 1460
$ws.Range("P9").Value = 1460

# Row 10
$ws.Range("D10").Value = 44526
$ws.Range("J10").Value = 100
$ws.Range("O10").Value = "Provincia de Linares"

# Row 11
$ws.Range("D11").Value = 44496
$ws.Range("J11").Value = 550
$ws.Range("K11").Value = 1500
$ws.Range("L11").Value = 2000
$ws.Range("M11").Value = 1773
$ws.Range("N11").Value = "$/paquete"
$ws.Range("P11").Value = 1773
